# Update "想去人数" (interested-count) figures to the values captured at the
# later GitHub Pages build (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")      # sheet1
$wsLocal   = $wb.Worksheets.Item("本地生活")  # sheet3
$wsAll     = $wb.Worksheets.Item("全部类型")  # sheet4

# --- 展览 (sheet1) ---
$wsExhibit.Range("F4").Value  = 3313
$wsExhibit.Range("F5").Value  = 211
$wsExhibit.Range("F6").Value  = 4810
$wsExhibit.Range("F7").Value  = 464
$wsExhibit.Range("F9").Value  = 174
$wsExhibit.Range("F14").Value = 653
$wsExhibit.Range("F20").Value = 4744
$wsExhibit.Range("F21").Value = 15
$wsExhibit.Range("F25").Value = 15
$wsExhibit.Range("F31").Value = 91
$wsExhibit.Range("F36").Value = 781
$wsExhibit.Range("F37").Value = 822

# --- 本地生活 (sheet3) ---
$wsLocal.Range("F3").Value = 1090

# --- 全部类型 (sheet4) ---
$wsAll.Range("F4").Value  = 1090
$wsAll.Range("F8").Value  = 3313
$wsAll.Range("F9").Value  = 211
$wsAll.Range("F10").Value = 4810
$wsAll.Range("F11").Value = 464
$wsAll.Range("F13").Value = 174
$wsAll.Range("F18").Value = 653
$wsAll.Range("F25").Value = 4744
$wsAll.Range("F26").Value = 15
$wsAll.Range("F30").Value = 15
$wsAll.Range("F37").Value = 91
$wsAll.Range("F42").Value = 781
$wsAll.Range("F43").Value = 822
